# Updated symbol list on Sun Jan 15 09:52:03 UTC 2023 with GitHub Actions
# Applies refreshed Price (D) and Volume(1h) (E) values for the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'294.70"
$ws.Range("E2").Value = "'-4.04%"
$ws.Range("D3").Value = "'31.24"
$ws.Range("E3").Value = "'-1.77%"
$ws.Range("D4").Value = "'5.097"
$ws.Range("E4").Value = "'-3.57%"
$ws.Range("D5").Value = "'0.07355"
$ws.Range("E5").Value = "'0.49%"
$ws.Range("D6").Value = "'7.681"
$ws.Range("E6").Value = "'-2.01%"
$ws.Range("D7").Value = "'3.751"
$ws.Range("E7").Value = "'-0.18%"
$ws.Range("D8").Value = "'1.653"
$ws.Range("E8").Value = "'11.89%"
$ws.Range("D9").Value = "'0.9198"
$ws.Range("E9").Value = "'1.51%"
$ws.Range("D10").Value = "'0.1672"
$ws.Range("E10").Value = "'-0.45%"
$ws.Range("D11").Value = "'0.07079"
$ws.Range("E11").Value = "'-5.42%"
$ws.Range("D12").Value = "'0.07922"
$ws.Range("E12").Value = "'0.14%"
$ws.Range("D13").Value = "'0.02997"
$ws.Range("E13").Value = "'-1.16%"
$ws.Range("D14").Value = "'0.09881"
$ws.Range("E14").Value = "'-0.95%"
$ws.Range("D15").Value = "'0.001491"
$ws.Range("E15").Value = "'-1.26%"
$ws.Range("D16").Value = "'0.006166"
$ws.Range("E16").Value = "'-1.10%"
$ws.Range("D17").Value = "'3.449"
$ws.Range("E17").Value = "'-1.16%"
$ws.Range("E18").Value = "'-0.08%"
$ws.Range("E19").Value = "'-1.41%"
$ws.Range("E20").Value = "'0.26%"
$ws.Range("D21").Value = "'4.559"
$ws.Range("E21").Value = "'6.86%"
$ws.Range("D22").Value = "'0.04611"
$ws.Range("E22").Value = "'0.97%"
$ws.Range("D24").Value = "'0.001218"
$ws.Range("E24").Value = "'-1.02%"
$ws.Range("D25").Value = "'0.004416"
$ws.Range("E25").Value = "'-0.04%"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'-0.58%"
$ws.Range("D27").Value = "'0.0001876"
$ws.Range("E27").Value = "'6.53%"
$ws.Range("D39").Value = "'0.01688"
$ws.Range("E39").Value = "'3.20%"
$ws.Range("D40").Value = "'0.04399"
$ws.Range("E40").Value = "'-2.37%"
$ws.Range("D41").Value = "'0.007079"
$ws.Range("E41").Value = "'0.09%"
$ws.Range("D42").Value = "'0.1324"
$ws.Range("E42").Value = "'-1.36%"
$ws.Range("D43").Value = "'0.002109"
$ws.Range("E43").Value = "'-7.18%"
$ws.Range("D44").Value = "'0.01045"
$ws.Range("E44").Value = "'-26.82%"
$ws.Range("D45").Value = "'0.00005978"
$ws.Range("E45").Value = "'-1.74%"
$ws.Range("D46").Value = "'1.918"
$ws.Range("E46").Value = "'1.34%"
$ws.Range("D47").Value = "'0.01100"
$ws.Range("E47").Value = "'-16.33%"